# Updates cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 per the Apr 7 2023 03:12:45 UTC GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so values like "28.136.92" or "313.32" are stored
# as literal text (matching the source t="inlineStr" cells) rather than
# being auto-coerced into numbers/dates by Excel value parsing.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.136.92'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.50'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.32'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5026'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3825'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08534'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -7.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.64'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.263'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.876.79'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.210'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001098'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '91.18'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06637'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.10'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.090'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.174.74'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.21'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.273'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.596'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.089.03'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.68'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '156.24'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '126.28'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1054'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.048'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -5.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.635'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.606'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.657'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02454'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.45%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06531'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.239'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2175'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.234'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6374'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.55%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.35'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.884'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6029'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.07'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.297'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.676'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.996'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.217'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '120.81'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '80.64'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +2.20%  '
